$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = '''20935'
$ws.Cells.Item(3, 1).Value = '''1817'
$ws.Cells.Item(4, 1).Value = '''3823'
$ws.Cells.Item(5, 1).Value = '''2785'
$ws.Cells.Item(6, 1).Value = '''4188'
$ws.Cells.Item(7, 1).Value = '''1278'
$ws.Cells.Item(8, 1).Value = '''7630'
$ws.Cells.Item(9, 1).Value = '''6438'
$ws.Cells.Item(10, 1).Value = '''386'
$ws.Cells.Item(11, 1).Value = '''5082'
$ws.Cells.Item(12, 1).Value = '''54751'
$ws.Cells.Item(13, 1).Value = '''1181'
$ws.Cells.Item(14, 1).Value = '''4927'
$ws.Cells.Item(15, 1).Value = '''3455'
$ws.Cells.Item(16, 1).Value = '''2375'
$ws.Cells.Item(17, 1).Value = '''16208'
$ws.Cells.Item(18, 1).Value = '''12107'
$ws.Cells.Item(19, 1).Value = '''4715'
$ws.Cells.Item(20, 1).Value = '''492'
$ws.Cells.Item(21, 1).Value = '''8689'
$ws.Cells.Item(22, 1).Value = '''4287'
$ws.Cells.Item(23, 1).Value = '''4813'
$ws.Cells.Item(24, 1).Value = '''3020'
$ws.Cells.Item(26, 1).Value = '''4173'
$ws.Cells.Item(26, 3).Value = '龍少'
$ws.Cells.Item(27, 1).Value = '''3002'
$ws.Cells.Item(28, 1).Value = '''3842'
$ws.Cells.Item(29, 1).Value = '''687'
$ws.Cells.Item(29, 3).Value = '"Smok3y 1nOnly"'
$ws.Cells.Item(30, 1).Value = '''7485'
$ws.Cells.Item(31, 1).Value = '''2119'
$ws.Cells.Item(32, 1).Value = '''1267'
$ws.Cells.Item(33, 1).Value = '''5437'
$ws.Cells.Item(34, 1).Value = '''4116'
$ws.Cells.Item(35, 1).Value = '''1552'
$ws.Cells.Item(36, 1).Value = '''4279'
$ws.Cells.Item(37, 1).Value = '''13727'
$ws.Cells.Item(38, 1).Value = '''4357'
$ws.Cells.Item(39, 1).Value = '''5467'
$ws.Cells.Item(40, 1).Value = '''902'
$ws.Cells.Item(41, 1).Value = '''6968'
$ws.Cells.Item(42, 1).Value = '''692'
$ws.Cells.Item(43, 1).Value = '''5258'
$ws.Cells.Item(44, 1).Value = '''4690'
$ws.Cells.Item(45, 1).Value = '''2571'
$ws.Cells.Item(46, 1).Value = '''4043'
$ws.Cells.Item(47, 1).Value = '''48733'
$ws.Cells.Item(48, 1).Value = '''2259'
$ws.Cells.Item(49, 1).Value = '''12135'
$ws.Cells.Item(50, 1).Value = '''1077'
$ws.Cells.Item(51, 1).Value = '''46100'
$ws.Cells.Item(52, 1).Value = '''50373'
$ws.Cells.Item(53, 1).Value = '''41006'
$ws.Cells.Item(54, 1).Value = '''5463'
$ws.Cells.Item(55, 1).Value = '''19763'
$ws.Cells.Item(56, 1).Value = '''16971'
$ws.Cells.Item(57, 1).Value = '''8015'
$ws.Cells.Item(58, 1).Value = '''14929'
$ws.Cells.Item(59, 1).Value = '''21420'
$ws.Cells.Item(60, 1).Value = '''20804'
$ws.Cells.Item(61, 1).Value = '''9223'
$ws.Cells.Item(63, 1).Value = '''8277'
$ws.Cells.Item(64, 1).Value = '''32503'
$ws.Cells.Item(65, 1).Value = '''16480'
$ws.Cells.Item(67, 1).Value = '''10774'
$ws.Cells.Item(68, 1).Value = '''2720'
$ws.Cells.Item(69, 1).Value = '''16834'
$ws.Cells.Item(70, 1).Value = '''17249'
$ws.Cells.Item(71, 1).Value = '''16118'
$ws.Cells.Item(72, 1).Value = '''23533'
$ws.Cells.Item(73, 1).Value = '''18382'
$ws.Cells.Item(74, 1).Value = '''17900'
$ws.Cells.Item(75, 1).Value = '''21803'
$ws.Cells.Item(76, 1).Value = '''9510'
$ws.Cells.Item(77, 1).Value = '''10680'
$ws.Cells.Item(78, 1).Value = '''7891'
$ws.Cells.Item(79, 1).Value = '''68304'
$ws.Cells.Item(80, 1).Value = '''12561'
$ws.Cells.Item(81, 1).Value = '''10630'
$ws.Cells.Item(82, 1).Value = '''12609'
$ws.Cells.Item(83, 1).Value = '''16102'
$ws.Cells.Item(84, 1).Value = '''41016'
$ws.Cells.Item(85, 1).Value = '''29246'
$ws.Cells.Item(86, 1).Value = '''33454'
$ws.Cells.Item(87, 1).Value = '''13635'
$ws.Cells.Item(88, 1).Value = '''49083'
$ws.Cells.Item(91, 1).Value = '''54910'
$ws.Cells.Item(92, 1).Value = '''77338'
$ws.Cells.Item(94, 1).Value = '''43436'
$ws.Cells.Item(95, 1).Value = '''57868'
$ws.Cells.Item(96, 1).Value = '''50830'
$ws.Cells.Item(97, 1).Value = '''58653'
$ws.Cells.Item(98, 1).Value = '''80905'
$ws.Cells.Item(99, 1).Value = '''26832'
$ws.Cells.Item(100, 3).Value = '人山即是仙'
$ws.Cells.Item(102, 1).Value = '''25851'
$ws.Cells.Item(103, 1).Value = '''31969'
$ws.Cells.Item(104, 1).Value = '''18612'
$ws.Cells.Item(105, 1).Value = '''26159'
$ws.Cells.Item(106, 1).Value = '''75102'
$ws.Cells.Item(107, 1).Value = '''59711'
$ws.Cells.Item(108, 1).Value = '''55610'
$ws.Cells.Item(109, 1).Value = '''27777'
$ws.Cells.Item(110, 1).Value = '''32579'
$ws.Cells.Item(111, 1).Value = '''46156'
$ws.Cells.Item(112, 1).Value = '''44728'
$ws.Cells.Item(113, 1).Value = '''33757'
$ws.Cells.Item(116, 1).Value = '''39451'
$ws.Cells.Item(117, 1).Value = '''27755'
$ws.Cells.Item(119, 1).Value = '''25416'
$ws.Cells.Item(121, 1).Value = '''45261'
$ws.Cells.Item(122, 1).Value = '''60726'
$ws.Cells.Item(123, 1).Value = '''39616'
$ws.Cells.Item(124, 1).Value = '''39885'
$ws.Cells.Item(125, 1).Value = '''30851'
$ws.Cells.Item(127, 1).Value = '''66020'
$ws.Cells.Item(130, 1).Value = '''50516'
$ws.Cells.Item(131, 1).Value = '''46069'
$ws.Cells.Item(133, 1).Value = '''47384'
$ws.Cells.Item(136, 1).Value = '''47413'
$ws.Cells.Item(140, 1).Value = '''47037'
$ws.Cells.Item(146, 1).Value = '''77255'
$ws.Cells.Item(149, 1).Value = '''6700'
$ws.Cells.Item(153, 1).Value = '''46712'
$ws.Cells.Item(155, 1).Value = '''47857'
